# "dev of alliance region"
# - Rename the "allianceShop" sheet to "shop"
# - Move the active-tab / selection state from "shrine" to "shop":
#     shrine: was the selected tab with active cell E2 -> becomes unselected, active cell D8
#     shop (ex allianceShop): was unselected with active cell G5 -> becomes the selected tab, active cell E8

$wb = $excel.ActiveWorkbook

# Update "shrine": it was the tab shown/selected in the original file.
# Activate it first so its new selection is recorded, then move the
# selection to D8 (tabSelected flag will be cleared once another sheet
# is activated below).
$wsShrine = $wb.Worksheets.Item("shrine")
$wsShrine.Activate()
[void]$wsShrine.Range("D8").Select()

# Rename "allianceShop" -> "shop"
$wsShop = $wb.Worksheets.Item("allianceShop")
$wsShop.Name = "shop"

# Make "shop" the active/selected tab and move its selection to E8.
$wsShop.Activate()
[void]$wsShop.Range("E8").Select()
